$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.033.27"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "3.065.75"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'517.21"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "'142.11"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.437"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("D10").Value = "'0.108"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "'0.376"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").Value = "3.585.89"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").Value = "'26.31"
$ws.Range("E14").Value = "  +3.78%  "
$ws.Range("D15").Value = "'0.0000164"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "58.012.60"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").Value = "3.059.00"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").Value = "'12.86"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").Value = "'8.07"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "'332.88"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'0.502"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "'65.48"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'0.171"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "0.0₃0905"
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").Value = "'6.46"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "'7.27"
$ws.Range("E29").Value = "  +7.07%  "
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("E31").Value = "  +4.05%  "
$ws.Range("D32").Value = "'20.72"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").Value = "'154.45"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "'4.54"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("D35").Value = "'6.02"
$ws.Range("E35").Value = "  +3.73%  "
$ws.Range("D36").Value = "'26.91"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("E37").Value = "  +3.55%  "
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("D39").Value = "3.104.42"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("D41").Value = "'36.63"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'0.658"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "2.293.37"
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("E45").Value = "  +5.94%  "
$ws.Range("D46").Value = "'1.38"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").Value = "'20.68"
$ws.Range("E47").Value = "  +4.56%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'5.94"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'0.939"
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("D50").Value = "'0.731"
$ws.Range("E50").Value = "  +8.80%  "
$ws.Range("D51").Value = "'0.0878"
$ws.Range("E51").Value = "  +2.69%  "
